$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.04271373187048222, 0.04071648406533734, 3.537761648806719, 1133.036916526867, 1136.65810839161)
    3 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    4 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    5 = @(0.1190320826869504, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 12.72000643965764)
    6 = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.611132179096228)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
